$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row (right answer marking value)
$ws.Range("B11").Value = 5

# Update "Total" row (total correct marks) and the corr/total marks label
$ws.Range("B12").Value = 130
$ws.Range("E12").Value = "130/140"
